$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the R:U values between the given row pairs.
$rowPairs = @(
    @(10, 11),
    @(12, 13),
    @(43, 44)
)

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    $rangeA = $ws.Range("R$rowA`:U$rowA")
    $rangeB = $ws.Range("R$rowB`:U$rowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}
